$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'4.13%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'40.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'7.86%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.759"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'11.91%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08103"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.70%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.597"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.16%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.782"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'3.49%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.962"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.29%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-1.18%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9455"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.87%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1284"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.42%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1993"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'3.07%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'8.993"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'37.75%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09263"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'2.68%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.03529"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'5.43%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.09638"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.33%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.001336"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.37%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006256"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'3.56%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.373"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.03%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3565"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.49%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1433"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'10.13%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2414"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'5.37%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04423"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.58%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001259"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'5.40%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004386"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.76%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001142"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-13.73%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'0.97%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02450"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'3.65%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05311"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.57%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007473"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.86%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1436"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'3.33%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.008738"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'2.06%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002124"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'6.66%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01081"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'23.43%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006887"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'8.29%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'1.16%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003177"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'11.37%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.001705"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'1.81%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002106"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'1.16%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'1.16%"
$ws.Range("E51").Style = "Normal"
